$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatting (style s="3" for Symbol, s="4" for Date) already used
# in rows 2:5 down to the new rows, mirroring the sheet's existing per-row
# formats so the appended rows look identical to the originals. Rows 9-10
# are a deliberate gap and are left completely untouched.
$ws.Range("A5:B5").Copy($ws.Range("A6:B8"))
$ws.Range("A5:B5").Copy($ws.Range("A11:B16"))

# Rows 2-7: symbols refreshed, watch date moved to 45383 (2024-04-01)
$ws.Range("A2").Value = "VSTE"
$ws.Range("A3").Value = "MDIA"
$ws.Range("A4").Value = "SGRP"
$ws.Range("A5").Value = "XLO"
$ws.Range("A6").Value = "MESO"
$ws.Range("A7").Value = "CXAI"
$ws.Range("B2:B7").Value = 45383

# Row 8 left blank (keeps formatting only, like the source)
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = ""

# Rows 9-10 intentionally skipped (no data/formatting touched)

# Rows 11-16: new symbols, watch date 45384 (2024-04-02)
$ws.Range("A11").Value = "VVPR"
$ws.Range("A12").Value = "AINC"
$ws.Range("A13").Value = "XTIA"
$ws.Range("A14").Value = "ADTX"
$ws.Range("A15").Value = "PIK"
$ws.Range("A16").Value = "LIFW"
$ws.Range("B11:B16").Value = 45384

# Match the updated selection shown in the diff
$ws.Range("B22").Select()
